$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2, D2, E2 and H2 hold plain text (a long numeric-looking reference id, a
# zero-padded code, an ISO date string and a free-text message) that must
# stay text rather than being auto-coerced into a number/date. Force text
# mode while assigning, then clear the format back so the cell keeps its
# original (unstyled) appearance, matching the source which carries no
# explicit style on those cells.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "72681594290101470231590"
$ws.Range("C2").ClearFormats()

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "143639"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2024-10-15"
$ws.Range("E2").ClearFormats()

$ws.Range("F2").Value = 92600

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "COUNTRY CODE INVALID FOR BUSINESS SERVICE ARRANGEMENT 4384001 AND ACCEPTANCE BRAND ID CODE 00000002 D0043 S06 DMC. INTERCHANGE RATE DESIGNATOR AND PROCESSING CODE/REVERSAL INDICATOR COMBINATION INVALID FOR 00000002 P0158 S04 BUSINESS SERVICE ARRANGEMENT 2060001 AND ACCEPTANCE BRAND ID CODE DMC."
$ws.Range("H2").ClearFormats()

# Remove rows 3 and 4 entirely (bottom-up so row numbers stay stable), which
# also shrinks the Table1 range + autoFilter + sheet dimension automatically.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()
